# Generate Report for Handoff
# Refresh the localization-status workbook with the newly generated
# handoff markdown file (new GUID) and xliff files (new content hash),
# along with the updated handoff timestamps.

$wb = $excel.ActiveWorkbook

$newGuid = "d6a5534b-124a-4535-8f3c-cfb62ed6460e"
$newHash = "a99fa954c34e681f7b56e8fb81b2a7a3bbc97707"

$newMdName   = "$newGuid.md"
$newMdPath   = "e2e\$newGuid.md"
$newZhXlf    = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlf    = "$newGuid.$newHash.de-de.xlf"

# The hyperlink target (commit blob URL) itself is unchanged by this edit -
# only the cell text / hyperlink display text is refreshed to the new file
# name, so reuse the existing external address for every rebuilt hyperlink.
$hlAddress = "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/3c39e48ee0d6d36a8502b8dbb79f82890b6d1b6f/e2e/38d20164-93f8-4271-b79e-67decc507e68.md"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("B2").Value = $newMdPath
$wsOverview.Range("G2").Value = "2017-01-03 05:24:36"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hlAddress, "", "", $newMdPath)

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newMdName
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = "2017-01-03 05:24:25"

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hlAddress, "", "", $newMdName)

# --- de-de sheet ---
# NOTE: de-de!H2 ("Latest Handback DateTime") happened to share the exact
# same shared-string slot as Overview!G2 in the source file (both held
# "2017-01-03 05:24:05"). The diff leaves de-de's own <row> markup
# untouched and only rewrites that shared-string's text, so H2's
# displayed value moves in lockstep with Overview!G2 - reproduce that by
# writing the same new text here too.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newMdName
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = "2017-01-03 05:24:36"

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hlAddress, "", "", $newMdName)
